$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.64875866666667
$ws.Range("H2").Value = 31.946276
$ws.Range("I2").Value = 0.04005553530491663
$ws.Range("J2").Value = 0.04005553530491663
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.106744
$ws.Range("N2").Value = 9.320231999999999
$ws.Range("O2").Value = 0.1418337845295607
$ws.Range("P2").Value = 0.1418337845295607
$ws.Range("Q2").Value = 33.08296709511466
$ws.Range("R2").Value = 297.746703856032
$ws.Range("S2").Value = 0.005681228163653755
$ws.Range("T2").Value = 0.005681228163653755

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.64875866666667
$ws.Range("H3").Value = 31.946276
$ws.Range("I3").Value = 0.04005553530491663
$ws.Range("J3").Value = 0.04005553530491663
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.522686333333333
$ws.Range("N3").Value = 4.568059
$ws.Range("O3").Value = 0.0695159837141737
$ws.Range("P3").Value = 0.06951598371417368
$ws.Range("Q3").Value = 16.21471928869822
$ws.Range("R3").Value = 145.932473598284
$ws.Range("S3").Value = 0.002784499939919094
$ws.Range("T3").Value = 0.002784499939919093

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.64875866666667
$ws.Range("H4").Value = 31.946276
$ws.Range("I4").Value = 0.04005553530491663
$ws.Range("J4").Value = 0.04005553530491663
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 17.274688
$ws.Range("N4").Value = 51.824064
$ws.Range("O4").Value = 0.7886502317562657
$ws.Range("P4").Value = 0.7886502317562656
$ws.Range("Q4").Value = 183.9539835539627
$ws.Range("R4").Value = 1655.585851985664
$ws.Range("S4").Value = 0.03158980720134378
$ws.Range("T4").Value = 0.03158980720134378

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 239.3334093333333
$ws.Range("H5").Value = 718.000228
$ws.Range("I5").Value = 0.9002577790786066
$ws.Range("J5").Value = 0.9002577790786066
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.106744
$ws.Range("N5").Value = 9.320231999999999
$ws.Range("O5").Value = 0.1418337845295607
$ws.Range("P5").Value = 0.1418337845295607
$ws.Range("Q5").Value = 743.5476334458772
$ws.Range("R5").Value = 6691.928701012896
$ws.Range("S5").Value = 0.1276869678588959
$ws.Range("T5").Value = 0.1276869678588959

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 239.3334093333333
$ws.Range("H6").Value = 718.000228
$ws.Range("I6").Value = 0.9002577790786066
$ws.Range("J6").Value = 0.9002577790786066
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.522686333333333
$ws.Range("N6").Value = 4.568059
$ws.Range("O6").Value = 0.0695159837141737
$ws.Range("P6").Value = 0.06951598371417368
$ws.Range("Q6").Value = 364.4297115019391
$ws.Range("R6").Value = 3279.867403517452
$ws.Range("S6").Value = 0.0625823051089866
$ws.Range("T6").Value = 0.06258230510898659

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 239.3334093333333
$ws.Range("H7").Value = 718.000228
$ws.Range("I7").Value = 0.9002577790786066
$ws.Range("J7").Value = 0.9002577790786066
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 17.274688
$ws.Range("N7").Value = 51.824064
$ws.Range("O7").Value = 0.7886502317562657
$ws.Range("P7").Value = 0.7886502317562656
$ws.Range("Q7").Value = 4134.409974209621
$ws.Range("R7").Value = 37209.68976788659
$ws.Range("S7").Value = 0.7099885061107242
$ws.Range("T7").Value = 0.7099885061107241

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 15.86769733333333
$ws.Range("H8").Value = 47.603092
$ws.Range("I8").Value = 0.05968668561647669
$ws.Range("J8").Value = 0.0596866856164767
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.106744
$ws.Range("N8").Value = 9.320231999999999
$ws.Range("O8").Value = 0.1418337845295607
$ws.Range("P8").Value = 0.1418337845295607
$ws.Range("Q8").Value = 49.29687348414933
$ws.Range("R8").Value = 443.671861357344
$ws.Range("S8").Value = 0.008465588507010983
$ws.Range("T8").Value = 0.008465588507010983

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 15.86769733333333
$ws.Range("H9").Value = 47.603092
$ws.Range("I9").Value = 0.05968668561647669
$ws.Range("J9").Value = 0.0596866856164767
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.522686333333333
$ws.Range("N9").Value = 4.568059
$ws.Range("O9").Value = 0.0695159837141737
$ws.Range("P9").Value = 0.06951598371417368
$ws.Range("Q9").Value = 24.16152587093645
$ws.Range("R9").Value = 217.453732838428
$ws.Range("S9").Value = 0.004149178665267999
$ws.Range("T9").Value = 0.004149178665267999

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 15.86769733333333
$ws.Range("H10").Value = 47.603092
$ws.Range("I10").Value = 0.05968668561647669
$ws.Range("J10").Value = 0.0596866856164767
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 17.274688
$ws.Range("N10").Value = 51.824064
$ws.Range("O10").Value = 0.7886502317562657
$ws.Range("P10").Value = 0.7886502317562656
$ws.Range("Q10").Value = 274.1095207117654
$ws.Range("R10").Value = 2466.985686405888
$ws.Range("S10").Value = 0.04707191844419771
$ws.Range("T10").Value = 0.04707191844419771
